$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 for 넥스트바이오메디컬 (new IPO entry), pushing
# the existing rows 3-21 down to 4-22.
$ws.Rows.Item(3).EntireRow.Insert()

$ws.Range("A3").Value = "넥스트바이오메디컬"
$ws.Range("B3").Value = "2024.07.15~07.19"
$ws.Range("C3").Value = "24,000~29,000"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = 24000
$ws.Range("F3").Value = "한국투자증권"

# The old 넥스트바이오메디컬 row (now shifted down to row 12) is superseded
# by the new entry above, so remove it.
$ws.Rows.Item(12).EntireRow.Delete()

# 한국스팩15호 now has a confirmed offering amount.
$ws.Range("D18").Value = "2000"
